$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "auto evaluvation and unique id for students" ---------------------
# The email column is replaced by an auto-generated, unique studentId
# column; the mailto: hyperlinks that lived on the email values go away
# with it.

# 1) Drop the mailto: hyperlinks attached to column A.
$ws.Hyperlinks.Delete()

# 2) Rename the column and populate it with unique, auto-generated
#    student ids (JAV202500001 .. JAV202500005), one per existing row.
$ws.Range("A1").Value = "studentId"

$studentIds = @("JAV202500001", "JAV202500002", "JAV202500003", "JAV202500004", "JAV202500005")
for ($i = 0; $i -lt $studentIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $studentIds[$i]
}

# 3) The studentId cells no longer need the hyperlink look (underline /
#    theme color) - give them a clean, wrapped & vertically centered
#    layout instead. Build the combined format on a cell that is already
#    inside the used range (so the sheet's dimension/used-range doesn't
#    balloon out), then copy just the formatting onto the id column.
$fmtSource = $ws.Range("F6")
$fmtSource.VerticalAlignment = -4108
$fmtSource.WrapText = $true
$fmtSource.Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)
$fmtSource.ClearFormats()

# 4) The "Hyperlink" cell style is no longer used anywhere - remove it.
$wb.Styles("Hyperlink").Delete()

# 5) Leave the cursor where the author left it after the edit.
$ws.Range("C9").Select() | Out-Null
